$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TSLA")

$ws.Range("B2").Value = 17446000000.0
$ws.Range("G2").Value = 6268000000.0

$ws.Range("B4").Value = 3600000000.0
$ws.Range("G4").Value = 3552000000.0

$ws.Range("B5").Value = 1237000000.0

$ws.Range("B7").Value = 24844000000.0
$ws.Range("G7").Value = 10396000000.0

$ws.Range("B8").Value = 277000000.0

$ws.Range("B10").Value = 1836000000.0

$ws.Range("B11").Value = 986000000.0
$ws.Range("G11").Value = 7608000000.0

$ws.Range("B15").Value = 7320000000.0
$ws.Range("G15").Value = 3771000000.0

$ws.Range("B17").Value = 2128000000.0

$ws.Range("B20").Value = 5273000000.0
$ws.Range("G20").Value = 726000000.0

$ws.Range("B22").Value = 10383000000.0
$ws.Range("G22").Value = 11634000000.0

$ws.Range("B24").Value = 144000000.0

$ws.Range("B26").Value = 1809000000.0
$ws.Range("G26").Value = 2691000000.0

$ws.Range("B28").Value = 28507000000.0
$ws.Range("G28").Value = 26842000000.0

$ws.Range("B32").Value = 23017000000.0
$ws.Range("G32").Value = 7467000000.0

$ws.Range("B33").Value = 23017000000.0
$ws.Range("G33").Value = 7467000000.0

$ws.Range("B37").Value = -4935000000.0
$ws.Range("G37").Value = 7436000000.0

$ws.Range("B38").Value = 12511000000.0
$ws.Range("G38").Value = 13704000000.0
